# Apply the "Spokane, WA Metro Area-Bus" FAC Summary Report update.
# Updates Year 1 from 2002 to 2012, refreshes the underlying factor data and
# ridership-effect figures for the new base year, switches the %Diff /
# Ridership-Effect columns to true ratios (displayed via number format
# rather than manually multiplying by 100), and updates the sheet's
# selection/scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Year 1 label: 2002 -> 2012 (top summary block and table header)
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "2012"
$ws.Range("E7").Value = "2012"

# ---------------------------------------------------------------------
# Row 8: Vehicle Revenue Miles
# ---------------------------------------------------------------------
$ws.Range("E8").Value = 5313529
$ws.Range("H8").Value = 438002.0607800001

# ---------------------------------------------------------------------
# Row 9: Average Fare (2018$)
# ---------------------------------------------------------------------
$ws.Range("E9").Value = 0.8816617729999999
$ws.Range("H9").Value = 327242.76318

# ---------------------------------------------------------------------
# Row 10: Population + Employment
# ---------------------------------------------------------------------
$ws.Range("E10").Value = 685240.5
$ws.Range("H10").Value = 649029.03171

# ---------------------------------------------------------------------
# Row 11: % of Population in Transit Supportive Density
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 20.48535424
$ws.Range("H11").Value = -182082.4548076

# ---------------------------------------------------------------------
# Row 12: Average Gas Price (2018$)
# ---------------------------------------------------------------------
$ws.Range("E12").Value = 4.3491
$ws.Range("H12").Value = -358191.7471099999

# ---------------------------------------------------------------------
# Row 13: Median Per Capita (2018$)
# ---------------------------------------------------------------------
$ws.Range("E13").Value = 27410.23
$ws.Range("H13").Value = -355834.61387

# ---------------------------------------------------------------------
# Row 14: % of Households with 0 Vehicles
# ---------------------------------------------------------------------
$ws.Range("E14").Value = 7.84
$ws.Range("H14").Value = -79001.2513399

# ---------------------------------------------------------------------
# Row 15: % Working at Home
# ---------------------------------------------------------------------
$ws.Range("E15").Value = 5.8
$ws.Range("H15").Value = 4486.243349999993

# ---------------------------------------------------------------------
# Row 19: New Reporters - Absolute Difference now explicitly 0
# ---------------------------------------------------------------------
$ws.Range("H19").Value = 0

# ---------------------------------------------------------------------
# Row 20: Total Modeled Ridership
# ---------------------------------------------------------------------
$ws.Range("E20").Value = 10056176.3

# ---------------------------------------------------------------------
# Row 21: Total Observed Ridership
# ---------------------------------------------------------------------
$ws.Range("E21").Value = 11030806

# ---------------------------------------------------------------------
# Formulas for the %Diff / Ridership Effect columns: drop the manual
# "*100" since the columns are now formatted as percentages directly.
# ---------------------------------------------------------------------
$pctRows = 8,9,10,11,12,13,14,15,16,17,18,19,20,21
foreach ($r in $pctRows) {
    $ws.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
$effRows = 8,9,10,11,12,13,14,15,16,17,18,19
foreach ($r in $effRows) {
    $ws.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}
# I20 / I21 simply mirror G20 / G21 (unchanged formulas, kept for clarity)
$ws.Range("I20").Formula = "=G20"
$ws.Range("I21").Formula = "=G21"

# ---------------------------------------------------------------------
# Number formats: the value columns (E, F, H) now show two decimals and
# the ratio columns (G, I) are displayed as percentages.
# ---------------------------------------------------------------------
$ws.Range("E8:F21").NumberFormat = "#,##0.00"
$ws.Range("H8:H21").NumberFormat = "#,##0.00"
$ws.Range("G8:G21").NumberFormat = "0.00%"
$ws.Range("I8:I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# Sheet view: drop the frozen/scrolled topLeftCell and move the active
# selection from K20 to H21.
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("A1").Select()
$ws.Range("H21").Select()
